$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115-134 down to 116-135
$ws.Rows.Item(115).EntireRow.Insert()

# Populate the newly inserted row 115 with the new weekly price record
$ws.Range("A115").Value = 11
$ws.Range("B115").Value = "Vega Monumental Concepción"
$ws.Range("C115").Value = "Bíobío"
$ws.Range("D115").Value = 44783
$ws.Range("E115").Value = 8
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100102
$ws.Range("H115").Value = "Cítricos"
$ws.Range("I115").Value = 100102004
$ws.Range("J115").Value = "Mandarina"
$ws.Range("K115").Value = "Murcott"
$ws.Range("L115").Value = "Primera"
$ws.Range("M115").Value = 320
$ws.Range("N115").Value = 7000
$ws.Range("O115").Value = 7500
$ws.Range("P115").Value = 7188
$ws.Range("Q115").Value = "$/caja 18 kilos"
$ws.Range("R115").Value = "Región de O'Higgins"
$ws.Range("S115").Value = 399
$ws.Range("T115").Value = 18
